$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: RFC_version_5 (same values as row 2 / RFC_version_3)
$ws.Range("A6").Value = "RFC_version_5"
$ws.Range("B6").Value = 0.3237458193979933
$ws.Range("C6").Value = 0.07548013755084859
$ws.Range("D6").Value = 0.07889785715663596
$ws.Range("E6").Value = 0.07663321681942116
$ws.Range("F6").Value = 0.293888166125396
$ws.Range("G6").Value = 0.3237458193979933
$ws.Range("H6").Value = 0.3073092323655619

# Row 7: RFC_version_6 (same values as row 3 / RFC_version_4)
$ws.Range("A7").Value = "RFC_version_6"
$ws.Range("B7").Value = 0.2909090909090909
$ws.Range("C7").Value = 0.1685804701627487
$ws.Range("D7").Value = 0.1551775519166823
$ws.Range("E7").Value = 0.1577336755908184
$ws.Range("F7").Value = 0.2757537399309551
$ws.Range("G7").Value = 0.2909090909090909
$ws.Range("H7").Value = 0.281001881001881
